# Ports workbook cleanup:
#   The "latitude" column (C) stores small negative decimal numbers as text
#   (they use the quote-prefix / "@ as text" style so Excel doesn't reinterpret
#   them). Three of the five recurring latitude labels were typed with a
#   French-locale comma instead of a period ("-3,456" / "-4,567" / "-5,678"),
#   inconsistent with their sibling values ("-1.234", "-2.345") and with the
#   "longitude" column (D) which already uses periods for the same numbers.
#   Fix those three labels everywhere they occur so the whole column is
#   period-formatted, matching the rest of the sheet's convention.
#
#   This pattern repeats every 5 data rows (rows 2-161, 32 blocks of 5):
#     offset 3 (rows 4, 9, 14, ...)  -> -3.456
#     offset 4 (rows 5, 10, 15, ...) -> -4.567
#     offset 5 (rows 6, 11, 16, ...) -> -5.678

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
[void]$ws.Activate()

$latitudeCol = 3   # column C

$row = 4
while ($row -le 159) {
    # Leading apostrophe forces text entry (keeps the existing quote-prefix
    # "stored as text" style instead of Excel re-parsing it as a number).
    $ws.Cells.Item($row, $latitudeCol).Value = "'-3.456"
    $row += 5
}

$row = 5
while ($row -le 160) {
    $ws.Cells.Item($row, $latitudeCol).Value = "'-4.567"
    $row += 5
}

$row = 6
while ($row -le 161) {
    $ws.Cells.Item($row, $latitudeCol).Value = "'-5.678"
    $row += 5
}

# Carry over the view state left behind when the workbook was last saved:
# zoomed out a bit, scrolled down toward the bottom of the data, and the
# active cell moved off the original selection.
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 122
$win.ScrollColumn = 1
[void]$ws.Range("J157").Select()
